{"js": "// Q1 answer-sheet cleanup:\n//  1. The first paragraph (\"Q1. ____ Shalu Bhati ____\") loses the bold\n//     \"Shalu Bhati\" name (and its spell-check proofing marks), the two\n//     blank/underscore runs around it get shortened, and a `_GoBack`\n//     bookmark is left where the name used to be.\n//  2. Both inline pictures (image2.png then image1.png) are removed,\n//     leaving their paragraphs empty. The `_GoBack` bookmark that used\n//     to sit next to the first picture is removed from there (it has\n//     been relocated into the first paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Drop the `_GoBack` bookmark from its old spot next to image2.png;\n// it gets re-inserted into the first paragraph below. ---\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// --- Remove every inline picture in the document (image2.png and\n// image1.png), leaving their host paragraphs empty. ---\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].inlinePictures.load(\"items\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const pictures = paragraphs.items[i].inlinePictures.items;\n  for (let j = 0; j < pictures.length; j++) {\n    pictures[j].delete();\n  }\n}\nawait context.sync();\n\n// --- Rewrite the first paragraph: shrink the two underscore blanks,\n// drop the bolded \"Shalu Bhati\" (and its proofing marks), and put a\n// fresh `_GoBack` bookmark in its place. ---\nconst firstParagraph = body.paragraphs.getFirst();\n\nconst newParagraphOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>Q1.</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\">      _____________</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n            <w:r><w:t>_______________________</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nfirstParagraph.getRange(\"Whole\").insertOoxml(newParagraphOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Q1 answer-sheet cleanup:\n#  1. The first paragraph (\"Q1. ____ Shalu Bhati ____\") loses the bold\n#     \"Shalu Bhati\" name (and its spell-check proofing marks), the two\n#     blank/underscore runs around it get shortened, and a `_GoBack`\n#     bookmark is left where the name used to be.\n#  2. Both inline pictures (image2.png then image1.png) are removed,\n#     leaving their paragraphs empty. The `_GoBack` bookmark that used\n#     to sit next to the first picture is removed from there (it has\n#     been relocated into the first paragraph).\n\n$d = $word.ActiveDocument\n\n# --- Drop the `_GoBack` bookmark from its old spot next to image2.png;\n# it gets re-inserted into the first paragraph below. ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Remove every inline picture in the document (image2.png and\n# image1.png), leaving their host paragraphs empty. ---\nwhile ($d.InlineShapes.Count -gt 0) {\n    $d.InlineShapes.Item(1).Delete()\n}\n\n# --- Rewrite the first paragraph: shrink the two underscore blanks,\n# drop the bolded \"Shalu Bhati\" (and its proofing marks), and put a\n# fresh `_GoBack` bookmark in its place. ---\n$firstParagraph = $d.Paragraphs.Item(1).Range\n\n$newParagraphXml = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>Q1.</w:t></w:r><w:r><w:t xml:space=\"preserve\">      _____________</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t>_______________________</w:t></w:r></w:p>\n'@\n\n$firstParagraph.InsertXML($newParagraphXml)\n"}
